$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New tracker snapshot date (2025-09-13 -> serial 45913)
$newDate = 45913
$progress = 0.9420452352542067

$data = @(
    @{ Row = 32; GoalID = "G2"; GoalName = "Workout" },
    @{ Row = 33; GoalID = "G3"; GoalName = "Eat Healthy" },
    @{ Row = 34; GoalID = "G4"; GoalName = "Read Book" },
    @{ Row = 35; GoalID = "G5"; GoalName = "Investment Plan" },
    @{ Row = 36; GoalID = "G6"; GoalName = "Spend 10 Hours without phone" }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.GoalID
    $ws.Cells.Item($r, 2).Value = $entry.GoalName
    $ws.Cells.Item($r, 3).Value = $newDate
    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 4).Value = $progress
    $ws.Cells.Item($r, 5).Value = 0
    $ws.Cells.Item($r, 6).Value = -0.01
}
